$d = $word.ActiveDocument

# Locate the very end of the document body content (after the last
# paragraph's text, "Not many write misses happen in other steps as the
# data is loaded and then stored back into.").
$endPos = $d.Content.End

# Re-acquire a fresh Range object at that position (re-fetching rather than
# reusing a Paragraph.Range that was Collapse()d avoids this host's stale
# Range.Text / content-eating quirk when the next InsertXML call runs).
$insertionPoint = $d.Range($endPos, $endPos)

# Two blank paragraphs, then the "In total cache misses ..." paragraph
# (three runs), then the "50,342,287 is the value cachegrind gave."
# paragraph (with the mid-sentence spell-check proofErr wrap around
# "cachegrind", matching the existing convention used elsewhere in this
# document for the same word).
$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p/><w:p/><w:p><w:r><w:t xml:space="preserve">In total cache misses 16,777,216 + </w:t></w:r><w:r><w:t>33,554,432</w:t></w:r><w:r><w:t xml:space="preserve"> = 50,331,648 estimated total cache misses (for the ones I have accounted for)</w:t></w:r></w:p><w:p><w:r><w:t>50,342,287</w:t></w:r><w:r><w:t xml:space="preserve"> is the value </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>cachegrind</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> gave.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$insertionPoint.InsertXML($xml) | Out-Null
